$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title change (appears twice: the H1 heading near the top, and the
#    bold run near the bottom) -> use Replace All so both get updated.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Kanpai Banzai Free: Slot Game Review", $true, $false, $false,
    $false, $false, $true, 1, $false, "Play Kanpai Banzai for Free", 2
) | Out-Null

# ------------------------------------------------------------------
# 2) Insert a new bullet "Wild reels, multipliers, scatters, and
#    cascading reels" right before "Potential to win over 10,000
#    times your bet" in the "What we like" list.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Great graphics and dynamic elements*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $rng = $target.Range.Duplicate
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    $newPara = $target.Next()
    $newRng = $newPara.Range.Duplicate
    $newRng.Collapse(1)
    $newRng.InsertAfter("Wild reels, multipliers, scatters, and cascading reels")
}

# ------------------------------------------------------------------
# 3) Remove the "Bonus game called Tasty Fish" bullet entirely.
# ------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Bonus game called Tasty Fish*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 4) "What we don't like" bullet text swaps.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Risk of premature end of bonus game", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Bonus game can end prematurely with poisoned piece", 2
) | Out-Null

$d.Content.Find.Execute(
    "Limited number of paylines", $true, $false, $false,
    $false, $false, $true, 1, $false, "RTP rate could be higher", 2
) | Out-Null

# ------------------------------------------------------------------
# 5) Meta-description (italic) paragraph near the very end.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Find out everything you need to know about the Kanpai Banzai slot game. Play for free and enjoy the great graphics, bonus games and winning potential.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Kanpai Banzai and play this exciting slot game for free.", 2
) | Out-Null

Write-Host "Edit complete"
